# --- LOT2042.xlsx update -----------------------------------------------
# Inserts a new row for "Docentes responsaveis" professor name (row 13),
# which pushes the remaining rows down by one, and refreshes all the
# Portuguese/English text content to the corrected/expanded wording.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Private Use Area bullet character used as a separator inside the existing
# English "Syllabus:" text (row 17, column B/C) - must be preserved as-is.
$bullet = [char]0xf02d

# 1. Insert a new row at position 13 - shifts old rows 13..23 down to 14..24,
#    (row heights / styles of the shifted rows come along automatically).
$ws.Rows.Item(13).Insert()

# The insert leaves a stray empty, styled A13 cell (copied from A12 above) -
# the target layout has no value in column A of row 13, so drop it entirely.
$ws.Range("A13").Clear()

# Give B13/C13 the same look (wrap text / red text) as the other data rows
# by copying the formatting from the row right below (still correctly styled).
$ws.Range("B14").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 2. Write the final text for every populated cell in the sheet.
$ws.Range("B1").Value = "Ementa atual:"
$ws.Range("C1").Value = "Ementa modificada (dados modificados em vermelho):"

$ws.Range("B2").Value = "LOT2042"
$ws.Range("C2").Value = "LOT2042"

$ws.Range("A3").Value = "Nome:"
$ws.Range("B3").Value = " Processos Bioquímicos"
$ws.Range("C3").Value = " Processos Bioquímicos"

$ws.Range("A4").Value = "Name:"
$ws.Range("B4").Value = "Biochemical Processes"
$ws.Range("C4").Value = "Biochemical Processes"

$ws.Range("A5").Value = "Créditos-aula:"
$ws.Range("B5").Value = "4"
$ws.Range("C5").Value = "4"

$ws.Range("A6").Value = "Créditos-trabalho"
$ws.Range("B6").Value = "0"
$ws.Range("C6").Value = "0"

$ws.Range("A7").Value = "Carga horária:"
$ws.Range("B7").Value = "60 h"
$ws.Range("C7").Value = "60 h"

$ws.Range("A8").Value = "Ativação:"
$ws.Range("B8").Value = "01/01/2018"
$ws.Range("C8").Value = "01/01/2018"

$ws.Range("A9").Value = "Semestre ideal:"
$ws.Range("B9").Value = "EQD-9,EQN-10"
$ws.Range("C9").Value = "EQD-9,EQN-10"

$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "Transmitir aos alunos os conceitos básicos e fundamentais no estudo de processos biotecnológicos. - Demonstrar as principais etapas no desenvolvimento dos processos bioquímicos industriais e os principais fatores que influenciam no desenvolvimento e otimização destes processos. - Apresentar aos alunos uma visão das aplicações potenciais e estratégicas da biotecnologia moderna - Aprimorar o raciocínio e despertar o espírito crítico e a criatividade dos alunos"
$ws.Range("C10").Value = "Transmitir aos alunos os conceitos básicos e fundamentais no estudo de processos biotecnológicos. - Demonstrar as principais etapas no desenvolvimento dos processos bioquímicos industriais e os principais fatores que influenciam no desenvolvimento e otimização destes processos. - Apresentar aos alunos uma visão das aplicações potenciais e estratégicas da biotecnologia moderna - Aprimorar o raciocínio e despertar o espírito crítico e a criatividade dos alunos"

$ws.Range("A11").Value = "Objectives:"

$ws.Range("A12").Value = "Docentes responsáveis:"

$ws.Range("B13").Value = "1814052 - Silvio Silverio da Silva"
$ws.Range("C13").Value = "1814052 - Silvio Silverio da Silva"

$ws.Range("A14").Value = "Programa resumido:"
$ws.Range("B14").Value = "Introdução. Processos fermentativos e enzimáticos. Processos bioquímicos aplicados à indústria de alimentos. Processos biotecnológicos de importância industrial. Tecnologia enzimática em diferentes setores industriais. Fundamentos de engenharia de bioprocessos aplicados a processos bioquímicos. Variação de escala em bioprocessos. Introdução às técnicas de separação/purificação de produtos biotecnológicos."
$ws.Range("C14").Value = "Introdução. Processos fermentativos e enzimáticos. Processos bioquímicos aplicados à indústria de alimentos. Processos biotecnológicos de importância industrial. Tecnologia enzimática em diferentes setores industriais. Fundamentos de engenharia de bioprocessos aplicados a processos bioquímicos. Variação de escala em bioprocessos. Introdução às técnicas de separação/purificação de produtos biotecnológicos."

$ws.Range("A15").Value = "Short syllabus:"
$ws.Range("B15").Value = "Introduction. Fermentative and enzymatic processes. Biochemical processes applied to food industry. Biochemical processes of industrial Importance. Enzymatic technology in different industrial sectors. Fundamentals of bioprocess engineering applied to biochemical processes. Variation of scale in bioprocesses. Introduction to techniques of separation/purification of biotechnological products."
$ws.Range("C15").Value = "Introduction. Fermentative and enzymatic processes. Biochemical processes applied to food industry. Biochemical processes of industrial Importance. Enzymatic technology in different industrial sectors. Fundamentals of bioprocess engineering applied to biochemical processes. Variation of scale in bioprocesses. Introduction to techniques of separation/purification of biotechnological products."

$ws.Range("A16").Value = "Programa:"
$ws.Range("B16").Value = "Introdução: importância dos bioprocessos e aplicações industriais.Processos fermentativos e enzimáticos: tipos de processos fermentativos, matérias primas, obtenção de produtos.Processos bioquímicos aplicados à indústria de alimentos: processamento de alimentos, fases do processamento de produtos alimentícios, alterações bioquímicas em alimentos, oxidação de lipídeos, escurecimento enzimático e não enzimático, controles industriais das alterações bioquímicas.Processos biotecnológicos de importância industrial: descrição e estudo de casos de alguns processos biotecnológicos.Tecnologia enzimática em diferentes setores industriais. Fundamentos de engenharia de bioprocessos aplicados a processos bioquímicos: transferência de oxigênio e respiração microbiana: transferência de massa (transferência por convecção em sistema gás-líquido; respiração microbiana; transferência de O2 da bolha de gás para a célula); transferência de O2 em biorreator (efeitos dos aspectos do dimensionamento e operacionais do biorreator - bolhas, aeração, agitação e propriedades do meio, agentes antiespumantes, temperatura, pressão do gás e pressão parcial de oxigênio). Transferência de potência e oxigênio em biorreator agitado e aerado. Ampliação de escala em bioprocessos.Introdução às técnicas de separação/purificação de produtos biotecnológicos."
$ws.Range("C16").Value = "Introdução: importância dos bioprocessos e aplicações industriais.Processos fermentativos e enzimáticos: tipos de processos fermentativos, matérias primas, obtenção de produtos.Processos bioquímicos aplicados à indústria de alimentos: processamento de alimentos, fases do processamento de produtos alimentícios, alterações bioquímicas em alimentos, oxidação de lipídeos, escurecimento enzimático e não enzimático, controles industriais das alterações bioquímicas.Processos biotecnológicos de importância industrial: descrição e estudo de casos de alguns processos biotecnológicos.Tecnologia enzimática em diferentes setores industriais. Fundamentos de engenharia de bioprocessos aplicados a processos bioquímicos: transferência de oxigênio e respiração microbiana: transferência de massa (transferência por convecção em sistema gás-líquido; respiração microbiana; transferência de O2 da bolha de gás para a célula); transferência de O2 em biorreator (efeitos dos aspectos do dimensionamento e operacionais do biorreator - bolhas, aeração, agitação e propriedades do meio, agentes antiespumantes, temperatura, pressão do gás e pressão parcial de oxigênio). Transferência de potência e oxigênio em biorreator agitado e aerado. Ampliação de escala em bioprocessos.Introdução às técnicas de separação/purificação de produtos biotecnológicos."

$ws.Range("A17").Value = "Syllabus:"
$ws.Range("B17").Value = $bullet + "Introduction: importance of bioprocesses and industrial applications." + $bullet + "Fermentative and enzymatic processes: types of fermentative processes, raw materials, products obtainment." + $bullet + "Biochemical processes applied to food industry: food processing, phases of food products processing, biochemical alterations in food, lipids oxidation, enzymatic and not enzymatic darkness, industrial controls of biochemical alterations." + $bullet + "Biotechnological processes of industrial importance: description and study of cases of some biotechnological processes." + $bullet + "Enzymatic technology in different industrial sectors. Fundamentals of bioprocess engineering applied to biochemical processes: transfer of oxygen and microbial respiration: transfer of mass (convection transfer in gas-liquid system, microbial respiration, transfer of oxygen from the gas bubble to the cell); Transfer of O2 in the bioreactor (Effects of sizing and operating aspects of the bioreactor - bubbles, aeration, agitation and medium properties, antifoaming agents, temperature, gas pressure and oxygen partial pressure). Transfer of power and oxygen in agitated and aerated bioreactor. Scale variation in bioprocesses.Introduction to separation/purification of biotechnological products"
$ws.Range("C17").Value = $bullet + "Introduction: importance of bioprocesses and industrial applications." + $bullet + "Fermentative and enzymatic processes: types of fermentative processes, raw materials, products obtainment." + $bullet + "Biochemical processes applied to food industry: food processing, phases of food products processing, biochemical alterations in food, lipids oxidation, enzymatic and not enzymatic darkness, industrial controls of biochemical alterations." + $bullet + "Biotechnological processes of industrial importance: description and study of cases of some biotechnological processes." + $bullet + "Enzymatic technology in different industrial sectors. Fundamentals of bioprocess engineering applied to biochemical processes: transfer of oxygen and microbial respiration: transfer of mass (convection transfer in gas-liquid system, microbial respiration, transfer of oxygen from the gas bubble to the cell); Transfer of O2 in the bioreactor (Effects of sizing and operating aspects of the bioreactor - bubbles, aeration, agitation and medium properties, antifoaming agents, temperature, gas pressure and oxygen partial pressure). Transfer of power and oxygen in agitated and aerated bioreactor. Scale variation in bioprocesses.Introduction to separation/purification of biotechnological products"

$ws.Range("A18").Value = "Avaliação:"

$ws.Range("A19").Value = "Método:"
$ws.Range("B19").Value = "Os alunos serão avaliados formalmente por uma prova teórica (P) e trabalhos (T). A ponderação das notas será de 70% para a prova teórica (P) e 30% para a média aritmética das notas dos trabalhos (T), ou seja: Média Final do período letivo normal (MF) = (0,7xP +0,3xT)."
$ws.Range("C19").Value = "Os alunos serão avaliados formalmente por uma prova teórica (P) e trabalhos (T). A ponderação das notas será de 70% para a prova teórica (P) e 30% para a média aritmética das notas dos trabalhos (T), ou seja: Média Final do período letivo normal (MF) = (0,7xP +0,3xT)."

$ws.Range("A20").Value = "Critério:"
$ws.Range("B20").Value = "Serão aprovados os alunos que obtiverem média do período letivo normal igual ou maior que 5."
$ws.Range("C20").Value = "Serão aprovados os alunos que obtiverem média do período letivo normal igual ou maior que 5."

$ws.Range("A21").Value = "Norma de recuperação:"
$ws.Range("B21").Value = "A recuperação será feita por meio de uma prova (PR) para alunos que tenham MF maior ou igual a 3,0 e menor do que 5,0 e pelo menos 70% de frequência. A nota de recuperação (NR) será a média simples entre a média final (MF) e a prova de recuperação (PR). Será considerado aprovado o aluno com NR maior ou igual a 5,0."
$ws.Range("C21").Value = "A recuperação será feita por meio de uma prova (PR) para alunos que tenham MF maior ou igual a 3,0 e menor do que 5,0 e pelo menos 70% de frequência. A nota de recuperação (NR) será a média simples entre a média final (MF) e a prova de recuperação (PR). Será considerado aprovado o aluno com NR maior ou igual a 5,0."

$ws.Range("A22").Value = "Bibliografia:"
$ws.Range("B22").Value = "1.Schmidell, W.; Lima, U. A.; Aquarone, E.; Borzani, W. Biotecnologia Industrial  Engenharia Bioquímica, vol. 2, São Paulo: Edgard Blücher, 2001. 2. Borzani, W.; Schmidell, W.; Lima, U. A.; Aquarone, E. Biotecnologia Industrial. Fundamentos Vol. 1. São Paulo: Ed. Edgard Blücher,  2001. 3. Manual of industrial microbiology and biotechnology - Demain, A. L.; Solomon, N.A. Eds.Washington, American Society for Microbiology, 1986. 4. Fermentation and Enzyme Technology - Wang, D.C. et al. New York, Wiley-Interscience, 1979. 5. Princípios de Tecnologia de Alimentos - Gava, A.J. São Paulo, Nobel, 1983. 6. LIMA , U. A et al. Biotecnología Industrial, Biotecnologia na produção de alimentos - Série Biotecnología, vol4. Ed. Edgard Blucher,Ltda , 2001. 7. Tecnologia de Alimentos - José Evangelista -Livraria"
$ws.Range("C22").Value = "1.Schmidell, W.; Lima, U. A.; Aquarone, E.; Borzani, W. Biotecnologia Industrial  Engenharia Bioquímica, vol. 2, São Paulo: Edgard Blücher, 2001. 2. Borzani, W.; Schmidell, W.; Lima, U. A.; Aquarone, E. Biotecnologia Industrial. Fundamentos Vol. 1. São Paulo: Ed. Edgard Blücher,  2001. 3. Manual of industrial microbiology and biotechnology - Demain, A. L.; Solomon, N.A. Eds.Washington, American Society for Microbiology, 1986. 4. Fermentation and Enzyme Technology - Wang, D.C. et al. New York, Wiley-Interscience, 1979. 5. Princípios de Tecnologia de Alimentos - Gava, A.J. São Paulo, Nobel, 1983. 6. LIMA , U. A et al. Biotecnología Industrial, Biotecnologia na produção de alimentos - Série Biotecnología, vol4. Ed. Edgard Blucher,Ltda , 2001. 7. Tecnologia de Alimentos - José Evangelista -Livraria"

$ws.Range("A23").Value = "Requisitos:"

$ws.Range("B24").Value = "LOT2041 -  Engenharia Bioquímica  (Requisito fraco)`n"
$ws.Range("C24").Value = "LOT2041 -  Engenharia Bioquímica  (Requisito fraco)`n"

# 3. Column A previously spanned cols 1-2 (min="1" max="2") even though col 2 has
#    its own distinct <col> definition; narrow it back down to just column 1.
$ws.Columns.Item(1).ColumnWidth = 30.7109375

